$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"22.31000000000005"
$ws.Range("H2").Value = [double]"5.287495430383871e-08"
$ws.Range("I2").Value = [double]"5.287495430383871e-08"
$ws.Range("L2").Value = [double]"41.95283281619419"
$ws.Range("M2").Value = "[27.823188287332286, 56.08247734505609]"
$ws.Range("N2").Value = [double]"3.341595282702059e-07"
$ws.Range("O2").Value = [double]"3.341595282702059e-07"
$ws.Range("P2").Value = [double]"1.389973926813502"
$ws.Range("Q2").Value = "[1.0126054398958093, 1.7673424137311944]"
$ws.Range("R2").Value = [double]"2.45014497579632e-09"
$ws.Range("S2").Value = [double]"2.45014497579632e-09"
$ws.Range("T2").Value = [double]"54.40726102116941"
$ws.Range("U2").Value = "[46.20234169769346, 62.61218034464536]"
$ws.Range("V2").Value = [double]"0"
$ws.Range("W2").Value = [double]"0"
$ws.Range("X2").Value = [double]"17.37455455455459"
$ws.Range("Y2").Value = [double]"16.03461461461465"
$ws.Range("Z2").Value = [double]"18.71449449449453"

# Row 3
$ws.Range("F3").Value = [double]"22.31000000000005"
$ws.Range("H3").Value = [double]"1.884160927279055e-08"
$ws.Range("I3").Value = [double]"1.884160927279055e-08"
$ws.Range("L3").Value = [double]"44.95846734607001"
$ws.Range("M3").Value = "[29.264944541115362, 60.651990151024656]"
$ws.Range("N3").Value = [double]"6.846423881157193e-07"
$ws.Range("O3").Value = [double]"6.846423881157193e-07"
$ws.Range("P3").Value = [double]"1.918289808498272"
$ws.Range("Q3").Value = "[1.54092132158058, 2.295658295415965]"
$ws.Range("R3").Value = [double]"2.473576898864849e-13"
$ws.Range("S3").Value = [double]"2.473576898864849e-13"
$ws.Range("T3").Value = [double]"56.7481580309257"
$ws.Range("U3").Value = "[48.17622854209695, 65.32008751975445]"
$ws.Range("V3").Value = [double]"0"
$ws.Range("W3").Value = [double]"0"
$ws.Range("X3").Value = [double]"15.49863863863867"
$ws.Range("Y3").Value = [double]"14.15869869869873"
$ws.Range("Z3").Value = [double]"16.83857857857861"

# Row 4
$ws.Range("F4").Value = [double]"22.31000000000005"
$ws.Range("H4").Value = [double]"1.043621067342571e-08"
$ws.Range("I4").Value = [double]"1.043621067342571e-08"
$ws.Range("L4").Value = [double]"49.12385932565423"
$ws.Range("M4").Value = "[31.107145618333746, 67.14057303297471]"
$ws.Range("N4").Value = [double]"1.762882257061449e-06"
$ws.Range("O4").Value = [double]"1.762882257061449e-06"
$ws.Range("P4").Value = [double]"2.371131992799504"
$ws.Range("Q4").Value = "[2.018921405009657, 2.72334258058935]"
$ws.Range("T4").Value = [double]"59.23194423314374"
$ws.Range("U4").Value = "[49.84213865472876, 68.62174981155871]"
$ws.Range("V4").Value = [double]"2.220446049250313e-16"
$ws.Range("W4").Value = [double]"2.220446049250313e-16"
$ws.Range("X4").Value = [double]"13.89071071071074"
$ws.Range("Y4").Value = [double]"12.64010010010013"
$ws.Range("Z4").Value = [double]"15.14132132132135"

# Row 5
$ws.Range("F5").Value = [double]"23.2600000000002"
$ws.Range("H5").Value = [double]"4.405562173959154e-07"
$ws.Range("I5").Value = [double]"4.405562173959154e-07"
$ws.Range("L5").Value = [double]"38.30954567971337"
$ws.Range("M5").Value = "[22.949260335912612, 53.66983102351412]"
$ws.Range("N5").Value = [double]"8.5119282264845e-06"
$ws.Range("O5").Value = [double]"8.5119282264845e-06"
$ws.Range("P5").Value = [double]"2.773658378845042"
$ws.Range("Q5").Value = "[2.3208161945438106, 3.2265005631462738]"
$ws.Range("R5").Value = [double]"4.440892098500626e-16"
$ws.Range("S5").Value = [double]"4.440892098500626e-16"
$ws.Range("T5").Value = [double]"52.96753568207945"
$ws.Range("U5").Value = "[44.37741487159238, 61.55765649256652]"
$ws.Range("V5").Value = [double]"4.440892098500626e-16"
$ws.Range("W5").Value = [double]"4.440892098500626e-16"
$ws.Range("X5").Value = [double]"12.99207207207218"
$ws.Range("Y5").Value = [double]"11.31567567567577"
$ws.Range("Z5").Value = [double]"14.66846846846859"

# Row 6
$ws.Range("F6").Value = [double]"23.2600000000002"
$ws.Range("H6").Value = [double]"3.091910475760962e-06"
$ws.Range("I6").Value = [double]"3.091910475760962e-06"
$ws.Range("L6").Value = [double]"39.38304188033679"
$ws.Range("M6").Value = "[23.377525075990548, 55.388558684683034]"
$ws.Range("N6").Value = [double]"1.065309666303804e-05"
$ws.Range("O6").Value = [double]"1.065309666303804e-05"
$ws.Range("P6").Value = [double]"2.962342622303889"
$ws.Range("Q6").Value = "[2.4591846397469648, 3.465500604860814]"
$ws.Range("R6").Value = [double]"1.998401444325282e-15"
$ws.Range("S6").Value = [double]"1.998401444325282e-15"
$ws.Range("T6").Value = [double]"57.9082666705612"
$ws.Range("U6").Value = "[48.160994874114515, 67.65553846700789]"
$ws.Range("V6").Value = [double]"1.332267629550188e-15"
$ws.Range("W6").Value = [double]"1.332267629550188e-15"
$ws.Range("X6").Value = [double]"12.29357357357368"
$ws.Range("Y6").Value = [double]"10.430910910911"
$ws.Range("Z6").Value = [double]"14.15623623623636"

# Row 7
$ws.Range("F7").Value = [double]"23.2600000000002"
$ws.Range("H7").Value = [double]"8.817170438213395e-09"
$ws.Range("I7").Value = [double]"8.817170438213395e-09"
$ws.Range("L7").Value = [double]"47.930211018106"
$ws.Range("M7").Value = "[34.0682548502877, 61.7921671859243]"
$ws.Range("N7").Value = [double]"1.151795947507139e-08"
$ws.Range("O7").Value = [double]"1.151795947507139e-08"
$ws.Range("P7").Value = [double]"-3.132158441416851"
$ws.Range("Q7").Value = "[-3.4843690292066976, -2.7799478536270046]"
$ws.Range("T7").Value = [double]"55.70424624575899"
$ws.Range("U7").Value = "[46.577244091595226, 64.83124839992274]"
$ws.Range("V7").Value = [double]"6.661338147750939e-16"
$ws.Range("W7").Value = [double]"6.661338147750939e-16"
$ws.Range("X7").Value = [double]"11.59507507507517"
$ws.Range("Y7").Value = [double]"10.2912112112113"
$ws.Range("Z7").Value = [double]"12.89893893893905"
